# Add 4 fields to employee sheet: birth, birthDate, gender, currentJob
# Inserted as new columns G:J (existing basicSalary..info shift right to K:W)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 new columns before the old "basicSalary" column (G) ---
$ws.Range("G1:J1").EntireColumn.Insert()

# --- Header row (row 1) for the new columns, left -> right ---
$ws.Cells.Item(1, 7).Value2  = "birth"       # G1
$ws.Cells.Item(1, 8).Value2  = "birthDate"   # H1
$ws.Cells.Item(1, 9).Value2  = "gender"      # I1
$ws.Cells.Item(1, 10).Value2 = "currentJob"  # J1

# Give the new header cells the same bordered / bold / centered look as the
# rest of row 1 (matches the existing cellXfs used for the other headers).
$hdrRng = $ws.Range($ws.Cells.Item(1, 7), $ws.Cells.Item(1, 10))
$hdrRng.Font.Bold = $true
$hdrRng.Borders.LineStyle = 1
$hdrRng.HorizontalAlignment = -4108
$hdrRng.VerticalAlignment = -4160

# H1 additionally carries the date number format (same header look + date fmt)
$ws.Cells.Item(1, 8).NumberFormat = "yyyy\-mm\-dd;@"

# --- Per-row data for the new columns ---
$birth      = @("Nipal", "India", "Nipal", "India", "India", "Nipal", "Nipal", "Nipal", "Nipal", "Nipal")
$birthDate  = @(38900, 38726, 38847, 38467, 39373, 38778, 38646, 39242, 38973, 38310)
$gender     = @("male", "male", "male", "male", "male", "male", "male", "male", "male", "male")
$currentJob = @("Driver", "Driver", "Driver", "Driver", "Driver", "Driver", "Driver", "Driver", "Driver", "Driver")

# Fill column by column (currentJob, birth, birthDate, gender) so that new
# shared-string entries land in the same order as the authored workbook.
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 2, 10).Value2 = $currentJob[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 2, 7).Value2 = $birth[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 8).Value2 = $birthDate[$i]
    $ws.Cells.Item($r, 8).NumberFormat = "yyyy\-mm\-dd;@"
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 2, 9).Value2 = $gender[$i]
}

# --- Apply the same date format to the pre-existing "...ExpireDate" columns ---
# After the insert these are P (passportExpireDate), R (qidExpireDate),
# T (hcExpireDate) and V (licenseExpireDate); header + data rows both restyle.
$dateCols = @(16, 18, 20, 22)
foreach ($col in $dateCols) {
    $ws.Range($ws.Cells.Item(1, $col), $ws.Cells.Item(11, $col)).NumberFormat = "yyyy\-mm\-dd;@"
}

# --- Column widths (approx. best-fit after the new columns / reformatting) ---
function Set-ColWidth($col, $width) {
    $ws.Columns($col).ColumnWidth = $width - (5/6)
}
Set-ColWidth 1 10.83203125
Set-ColWidth 2 7.6640625
Set-ColWidth 3 6.5
Set-ColWidth 4 7.5
Set-ColWidth 5 13
Set-ColWidth 6 23.33203125
Set-ColWidth 7 8.83203125
Set-ColWidth 8 10.1640625
Set-ColWidth 9 8
Set-ColWidth 10 9.6640625
Set-ColWidth 11 9.5
Set-ColWidth 12 9.6640625
Set-ColWidth 13 8.1640625
Set-ColWidth 14 8
Set-ColWidth 15 9.83203125
Set-ColWidth 16 16.33203125
Set-ColWidth 17 12
Set-ColWidth 18 12.1640625
Set-ColWidth 19 9.33203125
Set-ColWidth 20 11.5
Set-ColWidth 22 15.1640625

# --- Selection / active cell ---
$ws.Range("I17").Select()
